# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(8, 9).Value = 'b'
$ws.Cells.Item(8, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(12, 9).Value = 'sd'
$ws.Cells.Item(12, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(13, 9).Value = 'sv'
$ws.Cells.Item(13, 10).Value = 'Statement-opinion'
$ws.Cells.Item(33, 9).Value = 'b'
$ws.Cells.Item(33, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(45, 9).Value = 'b'
$ws.Cells.Item(45, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(92, 9).Value = 'sv'
$ws.Cells.Item(92, 10).Value = 'Statement-opinion'
$ws.Cells.Item(104, 9).Value = 'aa'
$ws.Cells.Item(104, 10).Value = 'Agree/Accept'
$ws.Cells.Item(113, 9).Value = 'ba'
$ws.Cells.Item(113, 10).Value = 'Appreciation'
$ws.Cells.Item(118, 9).Value = 'sd'
$ws.Cells.Item(118, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(120, 9).Value = 'sd'
$ws.Cells.Item(120, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(123, 9).Value = 'sv'
$ws.Cells.Item(123, 10).Value = 'Statement-opinion'
$ws.Cells.Item(133, 9).Value = 'sd'
$ws.Cells.Item(133, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(151, 9).Value = 'b'
$ws.Cells.Item(151, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(164, 9).Value = 'aa'
$ws.Cells.Item(164, 10).Value = 'Agree/Accept'
$ws.Cells.Item(170, 9).Value = 'b'
$ws.Cells.Item(170, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(176, 9).Value = 'sd'
$ws.Cells.Item(176, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(184, 9).Value = 'b'
$ws.Cells.Item(184, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(198, 9).Value = 'ba'
$ws.Cells.Item(198, 10).Value = 'Appreciation'
$ws.Cells.Item(200, 9).Value = 'sv'
$ws.Cells.Item(200, 10).Value = 'Statement-opinion'
$ws.Cells.Item(204, 9).Value = 'b'
$ws.Cells.Item(204, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(212, 9).Value = 'sv'
$ws.Cells.Item(212, 10).Value = 'Statement-opinion'
$ws.Cells.Item(237, 9).Value = 'sv'
$ws.Cells.Item(237, 10).Value = 'Statement-opinion'
$ws.Cells.Item(239, 9).Value = '%'
$ws.Cells.Item(239, 10).Value = 'Uninterpretable'
$ws.Cells.Item(240, 9).Value = 'aa'
$ws.Cells.Item(240, 10).Value = 'Agree/Accept'
$ws.Cells.Item(243, 9).Value = 'sd'
$ws.Cells.Item(243, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(246, 9).Value = '%'
$ws.Cells.Item(246, 10).Value = 'Uninterpretable'
$ws.Cells.Item(247, 9).Value = 'sv'
$ws.Cells.Item(247, 10).Value = 'Statement-opinion'
$ws.Cells.Item(258, 9).Value = 'sd'
$ws.Cells.Item(258, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(261, 9).Value = 'sd'
$ws.Cells.Item(261, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(279, 9).Value = 'b'
$ws.Cells.Item(279, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(297, 9).Value = 'b'
$ws.Cells.Item(297, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(309, 9).Value = 'ba'
$ws.Cells.Item(309, 10).Value = 'Appreciation'
$ws.Cells.Item(328, 9).Value = 'b'
$ws.Cells.Item(328, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(345, 9).Value = 'sv'
$ws.Cells.Item(345, 10).Value = 'Statement-opinion'
$ws.Cells.Item(351, 9).Value = 'sv'
$ws.Cells.Item(351, 10).Value = 'Statement-opinion'
$ws.Cells.Item(354, 9).Value = 'sv'
$ws.Cells.Item(354, 10).Value = 'Statement-opinion'
$ws.Cells.Item(359, 9).Value = 'b'
$ws.Cells.Item(359, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(360, 9).Value = 'sv'
$ws.Cells.Item(360, 10).Value = 'Statement-opinion'
$ws.Cells.Item(366, 9).Value = 'sv'
$ws.Cells.Item(366, 10).Value = 'Statement-opinion'
$ws.Cells.Item(368, 9).Value = 'sv'
$ws.Cells.Item(368, 10).Value = 'Statement-opinion'
$ws.Cells.Item(382, 9).Value = 'ba'
$ws.Cells.Item(382, 10).Value = 'Appreciation'
